# Rename the comparison-report column headers (row 1) so the "old"/"new"
# suffixes become the concrete format-version identifiers used for this
# merge (FV2210 = left-hand/"old" input file, FV2304 = right-hand/"new"
# input file), then turn the header range into a proper Excel Table and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns 1-10: "<name>_FV2210" (was "<name>_old")
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2210"
}

# Column 11 ("diff") is unchanged.
$ws.Cells.Item(1, 11).Value = "diff"

# Columns 12-21: "<name>_FV2304" (was "<name>_new")
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2304"
}

# Turn A1:U55 into an Excel Table (adds xl/tables/table1.xml, the
# tableParts entry on the sheet and the relationship to it) using the
# freshly written header row as the column names.
$tableRange = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add(1, $tableRange, 1, 1)
$lo.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
